$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("1. MainInformation")
$ws2 = $wb.Worksheets.Item("2. Customer")
$ws3 = $wb.Worksheets.Item("3. Guarantor")

# Prime the shared-strings table so "-" lands at index 0 and "Unexecuted"
# lands at index 1 (first-use order determines the <sst> index), then fill
# the rest of sheet 1 ("1. MainInformation", used range B1:F2, row 1 is a
# taller 30pt header row).
$ws1.Range("B2").Value = "-"
$ws1.Range("B1").Value = "Unexecuted"
$ws1.Range("C1:F1").Value = "Unexecuted"
$ws1.Range("C2:F2").Value = "-"
$ws1.Rows.Item(1).RowHeight = 30

# ---------------------------------------------------------------
# Sheet 2 - "2. Customer": used range A1:H2.
# Column B keeps the workbook default ("Normal") style; columns C/D
# pick up the general-alignment style (matching A/E-H) instead of the
# column's own left-aligned default.
# ---------------------------------------------------------------
$ws2.Range("A1").Value = "-"
$ws2.Range("B1").Value = "Unexecuted"
$ws2.Range("C1:F1").Value = "Unexecuted"
$ws2.Range("G1:H1").Value = "-"

$ws2.Range("A2").Value = "-"
$ws2.Range("B2").Value = "-"
$ws2.Range("C2:F2").Value = "-"
$ws2.Range("G2:H2").Value = "-"

$ws2.Range("B1:B2").Style = "Normal"
$ws2.Range("C1:D2").HorizontalAlignment = 1

[void]$ws2.Range("A1:H2").Select()
[void]$ws3.Activate()

# ---------------------------------------------------------------
# Sheet 3 - "3. Guarantor": used range B1:F2 (columns B-D and E-F each
# inherit their column's own default style, so no overrides needed).
# ---------------------------------------------------------------
$ws3.Range("B1:D1").Value = "Unexecuted"
$ws3.Range("E1:F1").Value = "Unexecuted"
$ws3.Range("B2:D2").Value = "-"
$ws3.Range("E2:F2").Value = "-"
